$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updates to existing FAQ answers (rows 2-37) ---
$c8 = @'
Te comparto los pasos para que formes parte de esta aventura! 1. Realiza tu depósito o transferencia para confirmar tu cupo 2. Envía el comprobante de pago por este medio 3. En 48 horas recibirás un código y link para completar tu registro y listo!
'@
$ws.Cells.Item(8, 3).Value = $c8

$c10 = @'
Una vez que hayas enviado el comprobante de depósito o transferencia vas a recibir un código y link para completar tu registro.
'@
$ws.Cells.Item(10, 3).Value = $c10

$c22 = @'
Hola. El Altar Reto Trail se realizará el domingo 19 de Octubre, con punto de partida en la parroquia Quimiag, a 20 minutos de Riobamba.
Las distancias disponibles son: 5K, 10K, 25K, 40K y 70K.
Este evento está avalado por ASET como clasificatorio para la selección nacional de trail y cuenta con índice UTMB, válido para competencias internacionales.
'@
$ws.Cells.Item(22, 3).Value = $c22

$c23 = @'
Hola.La Ruta del Hielero se realizará este 4 de mayo, parte desde la parroquia de San Juan a 25 minutos de Riobamba.
Tenemos 5 distancias, 5k, 10k, 25k, 40k y 70k.
Esta carrera está calificada por la ASET, Asociación Ecuatoriana de Trail como evento clasificatorio para la selección nacional de trail, que representará al País en torneos internacionales.
Además cuenta con INDEX UTMB
'@
$ws.Cells.Item(23, 3).Value = $c23

$ws.Cells.Item(35, 3).Value = ""

$c37 = @'
Si, para completar tu inscripción envía tu comprobante de pago por Whatsapp y en 48 horas recibirás un código único  y link para completar tu inscripción
'@
$ws.Cells.Item(37, 3).Value = $c37

# --- Row 13 question + answer updated (kit delivery -> Altar specific) ---
$b13 = @'
Donde se entregarán los kits (Altar Reto Trail)?
'@
$ws.Cells.Item(13, 2).Value = $b13
$c13 = @'
¡Hora de equiparte para el desafío! La entrega de kits para Altar Reto Trail será en el concesionario Nissan Renault el 18 de Octubre de 11h00 a 18h00.
⚠️ Recuerda que no se entregan kits el día de la carrera bajo ninguna circunstancia.
🔄 Ten listo tu código QR de confirmación para agilizar el proceso.
Aquí te dejo la ubicación para que no te pierdas 👉 [link de Google Maps]
'@
$ws.Cells.Item(13, 3).Value = $c13

# --- New FAQ rows 38-42 ---
$ws.Cells.Item(38, 1).Value = "faq"
$b38 = @'
Donde se entregarán los kits (Ruta del Hielero)?
'@
$ws.Cells.Item(38, 2).Value = $b38
$c38 = @'
¡Hora de equiparte para el desafío! La entrega de kits para Ruta del Hielero será en el concesionario Nissan Renault el 3 de Mayo de 11h00 a 18h00.
⚠️ Recuerda que no se entregan kits el día de la carrera bajo ninguna circunstancia.
🔄 Ten listo tu código QR de confirmación para agilizar el proceso.
Aquí te dejo la ubicación para que no te pierdas 👉 [link de Google Maps]
'@
$ws.Cells.Item(38, 3).Value = $c38

$ws.Cells.Item(39, 1).Value = "faq"
$b39 = @'
Donde se entregarán los kits (Rio 21K)?
'@
$ws.Cells.Item(39, 2).Value = $b39
$c39 = @'
¡Hora de equiparte para el desafío! La entrega de kits para Rio 21K será en el concesionario Nissan Renault el 21 de Noviembre de 11h00 a 18h00.
⚠️ Recuerda que no se entregan kits el día de la carrera bajo ninguna circunstancia.
🔄 Ten listo tu código QR de confirmación para agilizar el proceso.
Aquí te dejo la ubicación para que no te pierdas 👉 [link de Google Maps]
'@
$ws.Cells.Item(39, 3).Value = $c39

$ws.Cells.Item(40, 1).Value = "faq"
$b40 = @'
Por qué se cambio de fecha la carrera ?
'@
$ws.Cells.Item(40, 2).Value = $b40
$c40 = @'
¡Hola! 👋
La fecha del Altar Reto Trail se cambió por motivos de fuerza mayor relacionados con la situación actual del país y el toque de queda en Chimborazo, que dificultaban la logística y seguridad del evento.
La nueva fecha es el 19 de octubre, elegida para afectar lo menos posible a nuestra comunidad de corredores 🏔️
'@
$ws.Cells.Item(40, 3).Value = $c40

$ws.Cells.Item(41, 1).Value = "faq"
$b41 = @'
Ya tengo el código para completar mi inscripción, qué sigue?
'@
$ws.Cells.Item(41, 2).Value = $b41
$c41 = @'
Por favor para completar tu inscripción y llenar tus datos finales, sigue el siguiente link : https://naftaecplus.com/carreras/
'@
$ws.Cells.Item(41, 3).Value = $c41

$ws.Cells.Item(42, 1).Value = "faq"
$b42 = @'
No me llega el código de confirmación, qué puedo hacer?
'@
$ws.Cells.Item(42, 2).Value = $b42
$c42 = @'
Si es que ya han pasado 48 horas y no recibes tu código para seguir tu inscripción, por favor comunicate al siguiente número +593 99 423 9837
'@
$ws.Cells.Item(42, 3).Value = $c42

